$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Cells.Item(19, 2).Value2 = 6815304
$ws.Cells.Item(19, 6).Value2 = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(19, 7).Value2 = 'Sutjeska Niksic'
$ws.Cells.Item(19, 9).Value2 = 0
$ws.Cells.Item(19, 10).Value2 = 'D'
$ws.Cells.Item(19, 11).Value2 = 5.5
$ws.Cells.Item(19, 12).Value2 = 3.75
$ws.Cells.Item(19, 13).Value2 = 1.5
$ws.Cells.Item(19, 14).Value2 = 3.6
$ws.Cells.Item(19, 15).Value2 = 3.2
$ws.Cells.Item(19, 16).Value2 = 1.909
$ws.Cells.Item(19, 17).Value2 = 0.5
$ws.Cells.Item(19, 18).Value2 = 1.825
$ws.Cells.Item(19, 19).Value2 = 1.975
$ws.Cells.Item(19, 21).Value2 = 1.875
$ws.Cells.Item(19, 22).Value2 = 1.925
$ws.Cells.Item(19, 24).Value2 = 2.2
$ws.Cells.Item(19, 25).Value2 = -1
$ws.Cells.Item(19, 26).Value2 = 0.825
$ws.Cells.Item(19, 27).Value2 = -1
$ws.Cells.Item(19, 28).Value2 = -1
$ws.Cells.Item(19, 29).Value2 = 0.925

# Row 20
$ws.Cells.Item(20, 2).Value2 = 6815422
$ws.Cells.Item(20, 6).Value2 = 'OFK Mladost DG'
$ws.Cells.Item(20, 7).Value2 = 'FK Decic Tuzi'
$ws.Cells.Item(20, 9).Value2 = 2
$ws.Cells.Item(20, 10).Value2 = 'A'
$ws.Cells.Item(20, 11).Value2 = 2.4
$ws.Cells.Item(20, 12).Value2 = 3
$ws.Cells.Item(20, 13).Value2 = 2.75
$ws.Cells.Item(20, 14).Value2 = 3.1
$ws.Cells.Item(20, 15).Value2 = 3
$ws.Cells.Item(20, 16).Value2 = 2.15
$ws.Cells.Item(20, 17).Value2 = 0.25
$ws.Cells.Item(20, 18).Value2 = 1.875
$ws.Cells.Item(20, 19).Value2 = 1.925
$ws.Cells.Item(20, 21).Value2 = 2.025
$ws.Cells.Item(20, 22).Value2 = 1.775
$ws.Cells.Item(20, 24).Value2 = -1
$ws.Cells.Item(20, 25).Value2 = 1.15
$ws.Cells.Item(20, 26).Value2 = -1
$ws.Cells.Item(20, 27).Value2 = 0.925
$ws.Cells.Item(20, 28).Value2 = -0.5
$ws.Cells.Item(20, 29).Value2 = 0.3875

# Row 30
$ws.Cells.Item(30, 2).Value2 = 6815311
$ws.Cells.Item(30, 6).Value2 = 'FK Mornar Bar'
$ws.Cells.Item(30, 7).Value2 = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(30, 8).Value2 = 4
$ws.Cells.Item(30, 9).Value2 = 3
$ws.Cells.Item(30, 10).Value2 = 'H'
$ws.Cells.Item(30, 11).Value2 = 1.833
$ws.Cells.Item(30, 12).Value2 = 3.1
$ws.Cells.Item(30, 13).Value2 = 4
$ws.Cells.Item(30, 14).Value2 = 2.25
$ws.Cells.Item(30, 15).Value2 = 2.9
$ws.Cells.Item(30, 16).Value2 = 3.1
$ws.Cells.Item(30, 17).Value2 = -0.25
$ws.Cells.Item(30, 18).Value2 = 1.975
$ws.Cells.Item(30, 19).Value2 = 1.825
$ws.Cells.Item(30, 20).Value2 = 2
$ws.Cells.Item(30, 21).Value2 = 1.9
$ws.Cells.Item(30, 22).Value2 = 1.9
$ws.Cells.Item(30, 23).Value2 = 1.25
$ws.Cells.Item(30, 25).Value2 = -1
$ws.Cells.Item(30, 26).Value2 = 0.9750000000000001
$ws.Cells.Item(30, 27).Value2 = -1
$ws.Cells.Item(30, 28).Value2 = 0.8999999999999999
$ws.Cells.Item(30, 29).Value2 = -1

# Row 31
$ws.Cells.Item(31, 2).Value2 = 6815315
$ws.Cells.Item(31, 6).Value2 = 'FK Decic Tuzi'
$ws.Cells.Item(31, 7).Value2 = 'FK Rudar Pljevlja'
$ws.Cells.Item(31, 8).Value2 = 0
$ws.Cells.Item(31, 9).Value2 = 1
$ws.Cells.Item(31, 10).Value2 = 'A'
$ws.Cells.Item(31, 11).Value2 = 1.615
$ws.Cells.Item(31, 12).Value2 = 3.5
$ws.Cells.Item(31, 13).Value2 = 4.75
$ws.Cells.Item(31, 14).Value2 = 1.4
$ws.Cells.Item(31, 15).Value2 = 3.8
$ws.Cells.Item(31, 16).Value2 = 6.5
$ws.Cells.Item(31, 17).Value2 = -1.25
$ws.Cells.Item(31, 18).Value2 = 2
$ws.Cells.Item(31, 19).Value2 = 1.8
$ws.Cells.Item(31, 20).Value2 = 2.5
$ws.Cells.Item(31, 21).Value2 = 1.95
$ws.Cells.Item(31, 22).Value2 = 1.85
$ws.Cells.Item(31, 23).Value2 = -1
$ws.Cells.Item(31, 25).Value2 = 5.5
$ws.Cells.Item(31, 26).Value2 = -1
$ws.Cells.Item(31, 27).Value2 = 0.8
$ws.Cells.Item(31, 28).Value2 = -1
$ws.Cells.Item(31, 29).Value2 = 0.8500000000000001

# Row 32
$ws.Cells.Item(32, 2).Value2 = 6815316
$ws.Cells.Item(32, 6).Value2 = 'FK Rudar Pljevlja'
$ws.Cells.Item(32, 7).Value2 = 'FK Mornar Bar'
$ws.Cells.Item(32, 8).Value2 = 0
$ws.Cells.Item(32, 9).Value2 = 0
$ws.Cells.Item(32, 11).Value2 = 2.875
$ws.Cells.Item(32, 13).Value2 = 2.3
$ws.Cells.Item(32, 14).Value2 = 2.6
$ws.Cells.Item(32, 15).Value2 = 2.9
$ws.Cells.Item(32, 16).Value2 = 2.6
$ws.Cells.Item(32, 18).Value2 = 1.9
$ws.Cells.Item(32, 19).Value2 = 1.9
$ws.Cells.Item(32, 20).Value2 = 2
$ws.Cells.Item(32, 21).Value2 = 1.875
$ws.Cells.Item(32, 22).Value2 = 1.925
$ws.Cells.Item(32, 24).Value2 = 1.9
$ws.Cells.Item(32, 28).Value2 = -1
$ws.Cells.Item(32, 29).Value2 = 0.925

# Row 33
$ws.Cells.Item(33, 2).Value2 = 6815319
$ws.Cells.Item(33, 6).Value2 = 'FK Arsenal'
$ws.Cells.Item(33, 7).Value2 = 'OFK Petrovac'
$ws.Cells.Item(33, 8).Value2 = 2
$ws.Cells.Item(33, 9).Value2 = 2
$ws.Cells.Item(33, 11).Value2 = 2.5
$ws.Cells.Item(33, 13).Value2 = 2.625
$ws.Cells.Item(33, 14).Value2 = 2.5
$ws.Cells.Item(33, 15).Value2 = 3
$ws.Cells.Item(33, 16).Value2 = 2.625
$ws.Cells.Item(33, 18).Value2 = 1.825
$ws.Cells.Item(33, 19).Value2 = 1.975
$ws.Cells.Item(33, 20).Value2 = 2.25
$ws.Cells.Item(33, 21).Value2 = 2
$ws.Cells.Item(33, 22).Value2 = 1.8
$ws.Cells.Item(33, 24).Value2 = 2
$ws.Cells.Item(33, 28).Value2 = 1
$ws.Cells.Item(33, 29).Value2 = -1

# Row 38
$ws.Cells.Item(38, 2).Value2 = 6815321
$ws.Cells.Item(38, 6).Value2 = 'OFK Petrovac'
$ws.Cells.Item(38, 7).Value2 = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(38, 8).Value2 = 1
$ws.Cells.Item(38, 10).Value2 = 'D'
$ws.Cells.Item(38, 11).Value2 = 1.8
$ws.Cells.Item(38, 12).Value2 = 3.4
$ws.Cells.Item(38, 13).Value2 = 3.75
$ws.Cells.Item(38, 14).Value2 = 1.6
$ws.Cells.Item(38, 15).Value2 = 3.5
$ws.Cells.Item(38, 16).Value2 = 4.75
$ws.Cells.Item(38, 17).Value2 = -0.75
$ws.Cells.Item(38, 18).Value2 = 1.825
$ws.Cells.Item(38, 19).Value2 = 1.975
$ws.Cells.Item(38, 20).Value2 = 2.5
$ws.Cells.Item(38, 21).Value2 = 1.95
$ws.Cells.Item(38, 22).Value2 = 1.75
$ws.Cells.Item(38, 23).Value2 = -1
$ws.Cells.Item(38, 24).Value2 = 2.5
$ws.Cells.Item(38, 26).Value2 = -1
$ws.Cells.Item(38, 27).Value2 = 0.9750000000000001
$ws.Cells.Item(38, 28).Value2 = -1
$ws.Cells.Item(38, 29).Value2 = 0.75

# Row 39
$ws.Cells.Item(39, 2).Value2 = 6815322
$ws.Cells.Item(39, 6).Value2 = 'OFK Mladost DG'
$ws.Cells.Item(39, 7).Value2 = 'FK Arsenal'
$ws.Cells.Item(39, 8).Value2 = 2
$ws.Cells.Item(39, 10).Value2 = 'H'
$ws.Cells.Item(39, 11).Value2 = 2.375
$ws.Cells.Item(39, 12).Value2 = 3
$ws.Cells.Item(39, 13).Value2 = 2.75
$ws.Cells.Item(39, 14).Value2 = 2.625
$ws.Cells.Item(39, 15).Value2 = 3
$ws.Cells.Item(39, 16).Value2 = 2.4
$ws.Cells.Item(39, 17).Value2 = 0
$ws.Cells.Item(39, 18).Value2 = 2
$ws.Cells.Item(39, 19).Value2 = 1.8
$ws.Cells.Item(39, 20).Value2 = 2
$ws.Cells.Item(39, 21).Value2 = 1.725
$ws.Cells.Item(39, 22).Value2 = 1.975
$ws.Cells.Item(39, 23).Value2 = 1.625
$ws.Cells.Item(39, 24).Value2 = -1
$ws.Cells.Item(39, 26).Value2 = 1
$ws.Cells.Item(39, 27).Value2 = -1
$ws.Cells.Item(39, 28).Value2 = 0.7250000000000001
$ws.Cells.Item(39, 29).Value2 = -1

# Row 48
$ws.Cells.Item(48, 2).Value2 = 6815331
$ws.Cells.Item(48, 6).Value2 = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(48, 7).Value2 = 'FK Decic Tuzi'
$ws.Cells.Item(48, 8).Value2 = 0
$ws.Cells.Item(48, 10).Value2 = 'A'
$ws.Cells.Item(48, 11).Value2 = 3
$ws.Cells.Item(48, 12).Value2 = 3.1
$ws.Cells.Item(48, 13).Value2 = 2.2
$ws.Cells.Item(48, 14).Value2 = 5.75
$ws.Cells.Item(48, 16).Value2 = 1.5
$ws.Cells.Item(48, 17).Value2 = 1
$ws.Cells.Item(48, 18).Value2 = 1.825
$ws.Cells.Item(48, 19).Value2 = 1.975
$ws.Cells.Item(48, 21).Value2 = 1.95
$ws.Cells.Item(48, 22).Value2 = 1.85
$ws.Cells.Item(48, 24).Value2 = -1
$ws.Cells.Item(48, 25).Value2 = 0.5
$ws.Cells.Item(48, 26).Value2 = 0
$ws.Cells.Item(48, 27).Value2 = -0
$ws.Cells.Item(48, 28).Value2 = -1
$ws.Cells.Item(48, 29).Value2 = 0.8500000000000001

# Row 49
$ws.Cells.Item(49, 2).Value2 = 6815333
$ws.Cells.Item(49, 6).Value2 = 'Sutjeska Niksic'
$ws.Cells.Item(49, 7).Value2 = 'FK Jezero'
$ws.Cells.Item(49, 8).Value2 = 1
$ws.Cells.Item(49, 10).Value2 = 'D'
$ws.Cells.Item(49, 11).Value2 = 1.5
$ws.Cells.Item(49, 12).Value2 = 3.6
$ws.Cells.Item(49, 13).Value2 = 6
$ws.Cells.Item(49, 14).Value2 = 1.5
$ws.Cells.Item(49, 16).Value2 = 5.75
$ws.Cells.Item(49, 17).Value2 = -1
$ws.Cells.Item(49, 18).Value2 = 1.9
$ws.Cells.Item(49, 19).Value2 = 1.9
$ws.Cells.Item(49, 21).Value2 = 1.85
$ws.Cells.Item(49, 22).Value2 = 1.95
$ws.Cells.Item(49, 24).Value2 = 2.6
$ws.Cells.Item(49, 25).Value2 = -1
$ws.Cells.Item(49, 26).Value2 = -1
$ws.Cells.Item(49, 27).Value2 = 0.8999999999999999
$ws.Cells.Item(49, 28).Value2 = -0.5
$ws.Cells.Item(49, 29).Value2 = 0.475

# Row 53
$ws.Cells.Item(53, 2).Value2 = 6815426
$ws.Cells.Item(53, 6).Value2 = 'FK Decic Tuzi'
$ws.Cells.Item(53, 7).Value2 = 'Buducnost Podgorica'
$ws.Cells.Item(53, 8).Value2 = 2
$ws.Cells.Item(53, 9).Value2 = 0
$ws.Cells.Item(53, 10).Value2 = 'H'
$ws.Cells.Item(53, 11).Value2 = 3
$ws.Cells.Item(53, 12).Value2 = 3
$ws.Cells.Item(53, 13).Value2 = 2.25
$ws.Cells.Item(53, 14).Value2 = 2.75
$ws.Cells.Item(53, 15).Value2 = 3
$ws.Cells.Item(53, 16).Value2 = 2.4
$ws.Cells.Item(53, 17).Value2 = 0.25
$ws.Cells.Item(53, 18).Value2 = 1.7
$ws.Cells.Item(53, 19).Value2 = 2.1
$ws.Cells.Item(53, 20).Value2 = 2.25
$ws.Cells.Item(53, 21).Value2 = 1.95
$ws.Cells.Item(53, 22).Value2 = 1.85
$ws.Cells.Item(53, 23).Value2 = 1.75
$ws.Cells.Item(53, 25).Value2 = -1
$ws.Cells.Item(53, 26).Value2 = 0.7
$ws.Cells.Item(53, 27).Value2 = -1
$ws.Cells.Item(53, 28).Value2 = -0.5
$ws.Cells.Item(53, 29).Value2 = 0.425

# Row 54
$ws.Cells.Item(54, 2).Value2 = 6815334
$ws.Cells.Item(54, 6).Value2 = 'Sutjeska Niksic'
$ws.Cells.Item(54, 7).Value2 = 'FK Mornar Bar'
$ws.Cells.Item(54, 8).Value2 = 0
$ws.Cells.Item(54, 10).Value2 = 'A'
$ws.Cells.Item(54, 11).Value2 = 1.444
$ws.Cells.Item(54, 12).Value2 = 4
$ws.Cells.Item(54, 13).Value2 = 6.5
$ws.Cells.Item(54, 14).Value2 = 1.444
$ws.Cells.Item(54, 15).Value2 = 4
$ws.Cells.Item(54, 16).Value2 = 6.5
$ws.Cells.Item(54, 17).Value2 = -1.25
$ws.Cells.Item(54, 18).Value2 = 2
$ws.Cells.Item(54, 19).Value2 = 1.8
$ws.Cells.Item(54, 20).Value2 = 2.5
$ws.Cells.Item(54, 21).Value2 = 2
$ws.Cells.Item(54, 22).Value2 = 1.8
$ws.Cells.Item(54, 24).Value2 = -1
$ws.Cells.Item(54, 25).Value2 = 5.5
$ws.Cells.Item(54, 26).Value2 = -1
$ws.Cells.Item(54, 27).Value2 = 0.8
$ws.Cells.Item(54, 28).Value2 = -1
$ws.Cells.Item(54, 29).Value2 = 0.8

# Row 55
$ws.Cells.Item(55, 2).Value2 = 7279987
$ws.Cells.Item(55, 6).Value2 = 'FK Jezero'
$ws.Cells.Item(55, 7).Value2 = 'FK Arsenal'
$ws.Cells.Item(55, 8).Value2 = 1
$ws.Cells.Item(55, 9).Value2 = 1
$ws.Cells.Item(55, 10).Value2 = 'D'
$ws.Cells.Item(55, 11).Value2 = 2.1
$ws.Cells.Item(55, 13).Value2 = 3.25
$ws.Cells.Item(55, 14).Value2 = 2.05
$ws.Cells.Item(55, 16).Value2 = 3.4
$ws.Cells.Item(55, 17).Value2 = -0.25
$ws.Cells.Item(55, 18).Value2 = 1.8
$ws.Cells.Item(55, 19).Value2 = 2
$ws.Cells.Item(55, 20).Value2 = 2
$ws.Cells.Item(55, 21).Value2 = 1.925
$ws.Cells.Item(55, 22).Value2 = 1.875
$ws.Cells.Item(55, 23).Value2 = -1
$ws.Cells.Item(55, 24).Value2 = 2
$ws.Cells.Item(55, 26).Value2 = -0.5
$ws.Cells.Item(55, 27).Value2 = 0.5
$ws.Cells.Item(55, 28).Value2 = 0
$ws.Cells.Item(55, 29).Value2 = -0

# Row 62
$ws.Cells.Item(62, 2).Value2 = 7366684
$ws.Cells.Item(62, 6).Value2 = 'FK Rudar Pljevlja'
$ws.Cells.Item(62, 7).Value2 = 'OFK Petrovac'
$ws.Cells.Item(62, 8).Value2 = 1
$ws.Cells.Item(62, 9).Value2 = 0
$ws.Cells.Item(62, 10).Value2 = 'H'
$ws.Cells.Item(62, 11).Value2 = 2.875
$ws.Cells.Item(62, 12).Value2 = 2.9
$ws.Cells.Item(62, 13).Value2 = 2.375
$ws.Cells.Item(62, 14).Value2 = 2.625
$ws.Cells.Item(62, 15).Value2 = 2.9
$ws.Cells.Item(62, 16).Value2 = 2.55
$ws.Cells.Item(62, 17).Value2 = 0
$ws.Cells.Item(62, 18).Value2 = 1.925
$ws.Cells.Item(62, 19).Value2 = 1.875
$ws.Cells.Item(62, 20).Value2 = 2.25
$ws.Cells.Item(62, 21).Value2 = 1.925
$ws.Cells.Item(62, 22).Value2 = 1.875
$ws.Cells.Item(62, 23).Value2 = 1.625
$ws.Cells.Item(62, 24).Value2 = -1
$ws.Cells.Item(62, 26).Value2 = 0.925
$ws.Cells.Item(62, 27).Value2 = -1
$ws.Cells.Item(62, 28).Value2 = -1
$ws.Cells.Item(62, 29).Value2 = 0.875

# Row 63
$ws.Cells.Item(63, 2).Value2 = 7366683
$ws.Cells.Item(63, 6).Value2 = 'FK Arsenal'
$ws.Cells.Item(63, 7).Value2 = 'FK Mornar Bar'
$ws.Cells.Item(63, 9).Value2 = 2
$ws.Cells.Item(63, 10).Value2 = 'D'
$ws.Cells.Item(63, 11).Value2 = 2.375
$ws.Cells.Item(63, 12).Value2 = 2.8
$ws.Cells.Item(63, 13).Value2 = 3
$ws.Cells.Item(63, 14).Value2 = 2.3
$ws.Cells.Item(63, 15).Value2 = 2.7
$ws.Cells.Item(63, 16).Value2 = 3.3
$ws.Cells.Item(63, 17).Value2 = -0.25
$ws.Cells.Item(63, 18).Value2 = 2
$ws.Cells.Item(63, 19).Value2 = 1.8
$ws.Cells.Item(63, 20).Value2 = 1.75
$ws.Cells.Item(63, 21).Value2 = 1.875
$ws.Cells.Item(63, 22).Value2 = 1.925
$ws.Cells.Item(63, 23).Value2 = -1
$ws.Cells.Item(63, 24).Value2 = 1.7
$ws.Cells.Item(63, 26).Value2 = -0.5
$ws.Cells.Item(63, 27).Value2 = 0.4
$ws.Cells.Item(63, 28).Value2 = 0.875
$ws.Cells.Item(63, 29).Value2 = -1

# Row 64
$ws.Cells.Item(64, 2).Value2 = 6815343
$ws.Cells.Item(64, 6).Value2 = 'Sutjeska Niksic'
$ws.Cells.Item(64, 7).Value2 = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(64, 8).Value2 = 2
$ws.Cells.Item(64, 11).Value2 = 1.333
$ws.Cells.Item(64, 12).Value2 = 4.2
$ws.Cells.Item(64, 13).Value2 = 8
$ws.Cells.Item(64, 14).Value2 = 1.333
$ws.Cells.Item(64, 15).Value2 = 4.2
$ws.Cells.Item(64, 16).Value2 = 8
$ws.Cells.Item(64, 17).Value2 = -1.5
$ws.Cells.Item(64, 18).Value2 = 1.975
$ws.Cells.Item(64, 19).Value2 = 1.825
$ws.Cells.Item(64, 20).Value2 = 2.75
$ws.Cells.Item(64, 21).Value2 = 1.9
$ws.Cells.Item(64, 22).Value2 = 1.9
$ws.Cells.Item(64, 23).Value2 = 0.333
$ws.Cells.Item(64, 26).Value2 = 0.9750000000000001
$ws.Cells.Item(64, 29).Value2 = 0.8999999999999999

# Row 71
$ws.Cells.Item(71, 2).Value2 = 6815351
$ws.Cells.Item(71, 6).Value2 = 'FK Arsenal'
$ws.Cells.Item(71, 7).Value2 = 'Buducnost Podgorica'
$ws.Cells.Item(71, 8).Value2 = 1
$ws.Cells.Item(71, 10).Value2 = 'D'
$ws.Cells.Item(71, 11).Value2 = 4.2
$ws.Cells.Item(71, 12).Value2 = 3.6
$ws.Cells.Item(71, 13).Value2 = 1.666
$ws.Cells.Item(71, 14).Value2 = 5.25
$ws.Cells.Item(71, 15).Value2 = 3.8
$ws.Cells.Item(71, 16).Value2 = 1.5
$ws.Cells.Item(71, 17).Value2 = 1
$ws.Cells.Item(71, 18).Value2 = 1.85
$ws.Cells.Item(71, 19).Value2 = 1.95
$ws.Cells.Item(71, 20).Value2 = 2.5
$ws.Cells.Item(71, 21).Value2 = 1.8
$ws.Cells.Item(71, 22).Value2 = 2
$ws.Cells.Item(71, 24).Value2 = 2.8
$ws.Cells.Item(71, 25).Value2 = -1
$ws.Cells.Item(71, 26).Value2 = 0.8500000000000001
$ws.Cells.Item(71, 27).Value2 = -1
$ws.Cells.Item(71, 29).Value2 = 1

# Row 72
$ws.Cells.Item(72, 2).Value2 = 6815354
$ws.Cells.Item(72, 6).Value2 = 'FK Rudar Pljevlja'
$ws.Cells.Item(72, 7).Value2 = 'FK Decic Tuzi'
$ws.Cells.Item(72, 8).Value2 = 0
$ws.Cells.Item(72, 10).Value2 = 'A'
$ws.Cells.Item(72, 11).Value2 = 4
$ws.Cells.Item(72, 12).Value2 = 3.2
$ws.Cells.Item(72, 13).Value2 = 1.8
$ws.Cells.Item(72, 14).Value2 = 3.8
$ws.Cells.Item(72, 15).Value2 = 3.25
$ws.Cells.Item(72, 16).Value2 = 1.8
$ws.Cells.Item(72, 17).Value2 = 0.5
$ws.Cells.Item(72, 18).Value2 = 1.95
$ws.Cells.Item(72, 19).Value2 = 1.85
$ws.Cells.Item(72, 20).Value2 = 2.25
$ws.Cells.Item(72, 21).Value2 = 2
$ws.Cells.Item(72, 22).Value2 = 1.8
$ws.Cells.Item(72, 24).Value2 = -1
$ws.Cells.Item(72, 25).Value2 = 0.8
$ws.Cells.Item(72, 26).Value2 = -1
$ws.Cells.Item(72, 27).Value2 = 0.8500000000000001
$ws.Cells.Item(72, 29).Value2 = 0.8

# Row 81
$ws.Cells.Item(81, 2).Value2 = 6815362
$ws.Cells.Item(81, 6).Value2 = 'Sutjeska Niksic'
$ws.Cells.Item(81, 7).Value2 = 'FK Decic Tuzi'
$ws.Cells.Item(81, 8).Value2 = 1
$ws.Cells.Item(81, 9).Value2 = 1
$ws.Cells.Item(81, 10).Value2 = 'D'
$ws.Cells.Item(81, 11).Value2 = 2.2
$ws.Cells.Item(81, 12).Value2 = 3
$ws.Cells.Item(81, 13).Value2 = 3.1
$ws.Cells.Item(81, 14).Value2 = 2.375
$ws.Cells.Item(81, 15).Value2 = 2.875
$ws.Cells.Item(81, 16).Value2 = 3
$ws.Cells.Item(81, 17).Value2 = -0.25
$ws.Cells.Item(81, 18).Value2 = 2.05
$ws.Cells.Item(81, 19).Value2 = 1.75
$ws.Cells.Item(81, 20).Value2 = 2
$ws.Cells.Item(81, 21).Value2 = 1.8
$ws.Cells.Item(81, 22).Value2 = 2
$ws.Cells.Item(81, 23).Value2 = -1
$ws.Cells.Item(81, 24).Value2 = 1.875
$ws.Cells.Item(81, 27).Value2 = 0.375
$ws.Cells.Item(81, 28).Value2 = 0
$ws.Cells.Item(81, 29).Value2 = -0

# Row 82
$ws.Cells.Item(82, 2).Value2 = 6815430
$ws.Cells.Item(82, 6).Value2 = 'Buducnost Podgorica'
$ws.Cells.Item(82, 7).Value2 = 'FK Mornar Bar'
$ws.Cells.Item(82, 8).Value2 = 4
$ws.Cells.Item(82, 9).Value2 = 3
$ws.Cells.Item(82, 10).Value2 = 'H'
$ws.Cells.Item(82, 11).Value2 = 1.444
$ws.Cells.Item(82, 12).Value2 = 3.75
$ws.Cells.Item(82, 13).Value2 = 6.5
$ws.Cells.Item(82, 14).Value2 = 1.4
$ws.Cells.Item(82, 15).Value2 = 4
$ws.Cells.Item(82, 16).Value2 = 7
$ws.Cells.Item(82, 17).Value2 = -1.25
$ws.Cells.Item(82, 18).Value2 = 1.875
$ws.Cells.Item(82, 19).Value2 = 1.925
$ws.Cells.Item(82, 20).Value2 = 2.5
$ws.Cells.Item(82, 21).Value2 = 1.775
$ws.Cells.Item(82, 22).Value2 = 1.925
$ws.Cells.Item(82, 23).Value2 = 0.3999999999999999
$ws.Cells.Item(82, 24).Value2 = -1
$ws.Cells.Item(82, 27).Value2 = 0.4625
$ws.Cells.Item(82, 28).Value2 = 0.7749999999999999
$ws.Cells.Item(82, 29).Value2 = -1

# Row 107
$ws.Cells.Item(107, 2).Value2 = 7890506
$ws.Cells.Item(107, 6).Value2 = 'FK Mornar Bar'
$ws.Cells.Item(107, 7).Value2 = 'FK Arsenal'
$ws.Cells.Item(107, 8).Value2 = 0
$ws.Cells.Item(107, 9).Value2 = 0
$ws.Cells.Item(107, 11).Value2 = 1.85
$ws.Cells.Item(107, 13).Value2 = 3.9
$ws.Cells.Item(107, 14).Value2 = 1.85
$ws.Cells.Item(107, 15).Value2 = 3.3
$ws.Cells.Item(107, 16).Value2 = 3.5
$ws.Cells.Item(107, 18).Value2 = 1.925
$ws.Cells.Item(107, 19).Value2 = 1.875
$ws.Cells.Item(107, 20).Value2 = 2
$ws.Cells.Item(107, 24).Value2 = 2.3
$ws.Cells.Item(107, 27).Value2 = 0.875
$ws.Cells.Item(107, 28).Value2 = -1
$ws.Cells.Item(107, 29).Value2 = 0.8500000000000001

# Row 108
$ws.Cells.Item(108, 2).Value2 = 7890508
$ws.Cells.Item(108, 6).Value2 = 'OFK Petrovac'
$ws.Cells.Item(108, 7).Value2 = 'FK Rudar Pljevlja'
$ws.Cells.Item(108, 8).Value2 = 1
$ws.Cells.Item(108, 9).Value2 = 1
$ws.Cells.Item(108, 11).Value2 = 1.75
$ws.Cells.Item(108, 13).Value2 = 4.5
$ws.Cells.Item(108, 14).Value2 = 1.8
$ws.Cells.Item(108, 15).Value2 = 3.2
$ws.Cells.Item(108, 16).Value2 = 4
$ws.Cells.Item(108, 18).Value2 = 1.875
$ws.Cells.Item(108, 19).Value2 = 1.925
$ws.Cells.Item(108, 20).Value2 = 2.25
$ws.Cells.Item(108, 24).Value2 = 2.2
$ws.Cells.Item(108, 27).Value2 = 0.925
$ws.Cells.Item(108, 28).Value2 = -0.5
$ws.Cells.Item(108, 29).Value2 = 0.425

# Row 126
$ws.Cells.Item(126, 2).Value2 = 6815401
$ws.Cells.Item(126, 6).Value2 = 'FK Decic Tuzi'
$ws.Cells.Item(126, 7).Value2 = 'Sutjeska Niksic'
$ws.Cells.Item(126, 9).Value2 = 0
$ws.Cells.Item(126, 10).Value2 = 'D'
$ws.Cells.Item(126, 11).Value2 = 2.55
$ws.Cells.Item(126, 13).Value2 = 2.6
$ws.Cells.Item(126, 14).Value2 = 2.1
$ws.Cells.Item(126, 15).Value2 = 3.1
$ws.Cells.Item(126, 16).Value2 = 3.3
$ws.Cells.Item(126, 17).Value2 = -0.25
$ws.Cells.Item(126, 18).Value2 = 1.825
$ws.Cells.Item(126, 19).Value2 = 1.975
$ws.Cells.Item(126, 20).Value2 = 2
$ws.Cells.Item(126, 21).Value2 = 1.925
$ws.Cells.Item(126, 22).Value2 = 1.875
$ws.Cells.Item(126, 24).Value2 = 2.1
$ws.Cells.Item(126, 25).Value2 = -1
$ws.Cells.Item(126, 26).Value2 = -0.5
$ws.Cells.Item(126, 27).Value2 = 0.4875
$ws.Cells.Item(126, 29).Value2 = 0.875

# Row 127
$ws.Cells.Item(127, 2).Value2 = 6815402
$ws.Cells.Item(127, 6).Value2 = 'FK Rudar Pljevlja'
$ws.Cells.Item(127, 7).Value2 = 'FK Jezero'
$ws.Cells.Item(127, 9).Value2 = 1
$ws.Cells.Item(127, 10).Value2 = 'A'
$ws.Cells.Item(127, 11).Value2 = 2.8
$ws.Cells.Item(127, 13).Value2 = 2.375
$ws.Cells.Item(127, 14).Value2 = 2.45
$ws.Cells.Item(127, 15).Value2 = 2.9
$ws.Cells.Item(127, 16).Value2 = 2.75
$ws.Cells.Item(127, 17).Value2 = 0
$ws.Cells.Item(127, 18).Value2 = 1.775
$ws.Cells.Item(127, 19).Value2 = 2.025
$ws.Cells.Item(127, 20).Value2 = 1.75
$ws.Cells.Item(127, 21).Value2 = 1.825
$ws.Cells.Item(127, 22).Value2 = 1.975
$ws.Cells.Item(127, 24).Value2 = -1
$ws.Cells.Item(127, 25).Value2 = 1.75
$ws.Cells.Item(127, 26).Value2 = -1
$ws.Cells.Item(127, 27).Value2 = 1.025
$ws.Cells.Item(127, 29).Value2 = 0.9750000000000001

# Row 130
$ws.Cells.Item(130, 2).Value2 = 6815403
$ws.Cells.Item(130, 6).Value2 = 'FK Jezero'
$ws.Cells.Item(130, 7).Value2 = 'FK Mornar Bar'
$ws.Cells.Item(130, 8).Value2 = 3
$ws.Cells.Item(130, 11).Value2 = 2.875
$ws.Cells.Item(130, 12).Value2 = 2.75
$ws.Cells.Item(130, 13).Value2 = 2.5
$ws.Cells.Item(130, 14).Value2 = 3.1
$ws.Cells.Item(130, 15).Value2 = 2.55
$ws.Cells.Item(130, 16).Value2 = 2.55
$ws.Cells.Item(130, 17).Value2 = 0
$ws.Cells.Item(130, 18).Value2 = 2.05
$ws.Cells.Item(130, 19).Value2 = 1.75
$ws.Cells.Item(130, 20).Value2 = 1.75
$ws.Cells.Item(130, 21).Value2 = 1.95
$ws.Cells.Item(130, 22).Value2 = 1.85
$ws.Cells.Item(130, 23).Value2 = 2.1
$ws.Cells.Item(130, 26).Value2 = 1.05
$ws.Cells.Item(130, 27).Value2 = -1
$ws.Cells.Item(130, 28).Value2 = 0.95

# Row 131
$ws.Cells.Item(131, 2).Value2 = 6815404
$ws.Cells.Item(131, 6).Value2 = 'Sutjeska Niksic'
$ws.Cells.Item(131, 7).Value2 = 'FK Rudar Pljevlja'
$ws.Cells.Item(131, 8).Value2 = 2
$ws.Cells.Item(131, 9).Value2 = 1
$ws.Cells.Item(131, 10).Value2 = 'H'
$ws.Cells.Item(131, 11).Value2 = 1.5
$ws.Cells.Item(131, 12).Value2 = 3.75
$ws.Cells.Item(131, 13).Value2 = 5.75
$ws.Cells.Item(131, 14).Value2 = 1.5
$ws.Cells.Item(131, 15).Value2 = 3.75
$ws.Cells.Item(131, 16).Value2 = 5.75
$ws.Cells.Item(131, 17).Value2 = -1
$ws.Cells.Item(131, 18).Value2 = 1.85
$ws.Cells.Item(131, 19).Value2 = 1.95
$ws.Cells.Item(131, 20).Value2 = 2.25
$ws.Cells.Item(131, 21).Value2 = 1.9
$ws.Cells.Item(131, 22).Value2 = 1.9
$ws.Cells.Item(131, 23).Value2 = 0.5
$ws.Cells.Item(131, 25).Value2 = -1
$ws.Cells.Item(131, 26).Value2 = 0
$ws.Cells.Item(131, 27).Value2 = -0
$ws.Cells.Item(131, 28).Value2 = 0.8999999999999999
$ws.Cells.Item(131, 29).Value2 = -1

# Row 132
$ws.Cells.Item(132, 2).Value2 = 6815405
$ws.Cells.Item(132, 6).Value2 = 'FK Arsenal'
$ws.Cells.Item(132, 7).Value2 = 'FK Decic Tuzi'
$ws.Cells.Item(132, 9).Value2 = 2
$ws.Cells.Item(132, 10).Value2 = 'A'
$ws.Cells.Item(132, 11).Value2 = 4
$ws.Cells.Item(132, 12).Value2 = 3.2
$ws.Cells.Item(132, 13).Value2 = 1.8
$ws.Cells.Item(132, 14).Value2 = 5
$ws.Cells.Item(132, 15).Value2 = 3.25
$ws.Cells.Item(132, 16).Value2 = 1.65
$ws.Cells.Item(132, 17).Value2 = 0.75
$ws.Cells.Item(132, 18).Value2 = 1.875
$ws.Cells.Item(132, 19).Value2 = 1.925
$ws.Cells.Item(132, 20).Value2 = 2
$ws.Cells.Item(132, 21).Value2 = 1.875
$ws.Cells.Item(132, 22).Value2 = 1.925
$ws.Cells.Item(132, 24).Value2 = -1
$ws.Cells.Item(132, 25).Value2 = 0.6499999999999999
$ws.Cells.Item(132, 26).Value2 = -1
$ws.Cells.Item(132, 27).Value2 = 0.925
$ws.Cells.Item(132, 28).Value2 = 0
$ws.Cells.Item(132, 29).Value2 = -0

# Row 133
$ws.Cells.Item(133, 2).Value2 = 6815406
$ws.Cells.Item(133, 6).Value2 = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(133, 7).Value2 = 'OFK Mladost DG'
$ws.Cells.Item(133, 8).Value2 = 0
$ws.Cells.Item(133, 9).Value2 = 0
$ws.Cells.Item(133, 10).Value2 = 'D'
$ws.Cells.Item(133, 11).Value2 = 2.25
$ws.Cells.Item(133, 12).Value2 = 3.3
$ws.Cells.Item(133, 13).Value2 = 2.7
$ws.Cells.Item(133, 14).Value2 = 2.05
$ws.Cells.Item(133, 15).Value2 = 3.4
$ws.Cells.Item(133, 16).Value2 = 3
$ws.Cells.Item(133, 17).Value2 = -0.25
$ws.Cells.Item(133, 18).Value2 = 1.825
$ws.Cells.Item(133, 19).Value2 = 1.975
$ws.Cells.Item(133, 20).Value2 = 2.5
$ws.Cells.Item(133, 21).Value2 = 2
$ws.Cells.Item(133, 22).Value2 = 1.8
$ws.Cells.Item(133, 23).Value2 = -1
$ws.Cells.Item(133, 24).Value2 = 2.4
$ws.Cells.Item(133, 26).Value2 = -0.5
$ws.Cells.Item(133, 27).Value2 = 0.4875
$ws.Cells.Item(133, 28).Value2 = -1
$ws.Cells.Item(133, 29).Value2 = 0.8

# Row 135
$ws.Cells.Item(135, 2).Value2 = 8043517
$ws.Cells.Item(135, 6).Value2 = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(135, 7).Value2 = 'FK Decic Tuzi'
$ws.Cells.Item(135, 8).Value2 = 0
$ws.Cells.Item(135, 10).Value2 = 'A'
$ws.Cells.Item(135, 11).Value2 = 5.5
$ws.Cells.Item(135, 12).Value2 = 3.2
$ws.Cells.Item(135, 13).Value2 = 1.615
$ws.Cells.Item(135, 14).Value2 = 6.5
$ws.Cells.Item(135, 15).Value2 = 3.4
$ws.Cells.Item(135, 16).Value2 = 1.533
$ws.Cells.Item(135, 17).Value2 = 1
$ws.Cells.Item(135, 18).Value2 = 1.85
$ws.Cells.Item(135, 19).Value2 = 1.95
$ws.Cells.Item(135, 20).Value2 = 2.25
$ws.Cells.Item(135, 21).Value2 = 2.025
$ws.Cells.Item(135, 22).Value2 = 1.775
$ws.Cells.Item(135, 23).Value2 = -1
$ws.Cells.Item(135, 25).Value2 = 0.5329999999999999
$ws.Cells.Item(135, 26).Value2 = -1
$ws.Cells.Item(135, 27).Value2 = 0.95
$ws.Cells.Item(135, 28).Value2 = -0.5
$ws.Cells.Item(135, 29).Value2 = 0.3875

# Row 136
$ws.Cells.Item(136, 2).Value2 = 8043518
$ws.Cells.Item(136, 6).Value2 = 'FK Arsenal'
$ws.Cells.Item(136, 7).Value2 = 'FK Rudar Pljevlja'
$ws.Cells.Item(136, 8).Value2 = 4
$ws.Cells.Item(136, 10).Value2 = 'H'
$ws.Cells.Item(136, 11).Value2 = 1.909
$ws.Cells.Item(136, 12).Value2 = 3
$ws.Cells.Item(136, 13).Value2 = 3.9
$ws.Cells.Item(136, 14).Value2 = 1.65
$ws.Cells.Item(136, 15).Value2 = 3.3
$ws.Cells.Item(136, 16).Value2 = 5
$ws.Cells.Item(136, 17).Value2 = -0.75
$ws.Cells.Item(136, 18).Value2 = 1.875
$ws.Cells.Item(136, 19).Value2 = 1.925
$ws.Cells.Item(136, 20).Value2 = 2
$ws.Cells.Item(136, 21).Value2 = 1.8
$ws.Cells.Item(136, 22).Value2 = 2
$ws.Cells.Item(136, 23).Value2 = 0.6499999999999999
$ws.Cells.Item(136, 25).Value2 = -1
$ws.Cells.Item(136, 26).Value2 = 0.875
$ws.Cells.Item(136, 27).Value2 = -1
$ws.Cells.Item(136, 28).Value2 = 0.8
$ws.Cells.Item(136, 29).Value2 = -1

# Row 140
$ws.Cells.Item(140, 2).Value2 = 8062093
$ws.Cells.Item(140, 6).Value2 = 'FK Jezero'
$ws.Cells.Item(140, 7).Value2 = 'FK Arsenal'
$ws.Cells.Item(140, 8).Value2 = 4
$ws.Cells.Item(140, 9).Value2 = 0
$ws.Cells.Item(140, 10).Value2 = 'H'
$ws.Cells.Item(140, 11).Value2 = 2.1
$ws.Cells.Item(140, 13).Value2 = 3.25
$ws.Cells.Item(140, 14).Value2 = 2.1
$ws.Cells.Item(140, 15).Value2 = 3
$ws.Cells.Item(140, 16).Value2 = 3.2
$ws.Cells.Item(140, 18).Value2 = 1.875
$ws.Cells.Item(140, 19).Value2 = 1.925
$ws.Cells.Item(140, 23).Value2 = 1.1
$ws.Cells.Item(140, 24).Value2 = -1
$ws.Cells.Item(140, 26).Value2 = 0.875
$ws.Cells.Item(140, 27).Value2 = -1
$ws.Cells.Item(140, 28).Value2 = 0.95
$ws.Cells.Item(140, 29).Value2 = -1

# Row 141
$ws.Cells.Item(141, 2).Value2 = 8062094
$ws.Cells.Item(141, 6).Value2 = 'FK Rudar Pljevlja'
$ws.Cells.Item(141, 7).Value2 = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(141, 8).Value2 = 1
$ws.Cells.Item(141, 9).Value2 = 1
$ws.Cells.Item(141, 10).Value2 = 'D'
$ws.Cells.Item(141, 11).Value2 = 2.25
$ws.Cells.Item(141, 12).Value2 = 3
$ws.Cells.Item(141, 13).Value2 = 3
$ws.Cells.Item(141, 14).Value2 = 2.25
$ws.Cells.Item(141, 15).Value2 = 3.1
$ws.Cells.Item(141, 16).Value2 = 2.875
$ws.Cells.Item(141, 17).Value2 = -0.25
$ws.Cells.Item(141, 18).Value2 = 2
$ws.Cells.Item(141, 19).Value2 = 1.8
$ws.Cells.Item(141, 20).Value2 = 2.25
$ws.Cells.Item(141, 21).Value2 = 1.95
$ws.Cells.Item(141, 22).Value2 = 1.85
$ws.Cells.Item(141, 23).Value2 = -1
$ws.Cells.Item(141, 24).Value2 = 2.1
$ws.Cells.Item(141, 26).Value2 = -0.5
$ws.Cells.Item(141, 27).Value2 = 0.4
$ws.Cells.Item(141, 28).Value2 = -0.5
$ws.Cells.Item(141, 29).Value2 = 0.425

# Row 142
$ws.Cells.Item(142, 2).Value2 = 8062092
$ws.Cells.Item(142, 6).Value2 = 'Sutjeska Niksic'
$ws.Cells.Item(142, 7).Value2 = 'FK Mornar Bar'
$ws.Cells.Item(142, 8).Value2 = 2
$ws.Cells.Item(142, 11).Value2 = 1.65
$ws.Cells.Item(142, 12).Value2 = 3.2
$ws.Cells.Item(142, 13).Value2 = 5
$ws.Cells.Item(142, 14).Value2 = 1.8
$ws.Cells.Item(142, 16).Value2 = 4.5
$ws.Cells.Item(142, 17).Value2 = -0.5
$ws.Cells.Item(142, 18).Value2 = 1.825
$ws.Cells.Item(142, 19).Value2 = 1.975
$ws.Cells.Item(142, 20).Value2 = 1.75
$ws.Cells.Item(142, 21).Value2 = 1.775
$ws.Cells.Item(142, 22).Value2 = 2.025
$ws.Cells.Item(142, 23).Value2 = 0.8
$ws.Cells.Item(142, 26).Value2 = 0.825
$ws.Cells.Item(142, 28).Value2 = 0.3875
$ws.Cells.Item(142, 29).Value2 = -0.5

Write-Host "Row swap edits applied"